{"js": "// Apply the text replacements described by the diff:\n// - update the date line\n// - update each \"NNN\u00d7N=\" multiplication prompt in the table\nconst replacements = [\n  [\"2024-02-10 Saturday\", \"2024-02-11 Sunday\"],\n  [\"670\u00d76=\", \"392\u00d79=\"],\n  [\"394\u00d77=\", \"174\u00d75=\"],\n  [\"268\u00d73=\", \"424\u00d73=\"],\n  [\"704\u00d76=\", \"802\u00d76=\"],\n  [\"468\u00d75=\", \"514\u00d76=\"],\n  [\"989\u00d77=\", \"904\u00d72=\"],\n  [\"579\u00d79=\", \"518\u00d79=\"],\n  [\"946\u00d75=\", \"139\u00d75=\"],\n  [\"882\u00d75=\", \"431\u00d72=\"],\n  [\"652\u00d73=\", \"299\u00d72=\"],\n  [\"147\u00d76=\", \"333\u00d74=\"],\n  [\"939\u00d78=\", \"422\u00d72=\"],\n  [\"680\u00d72=\", \"115\u00d79=\"],\n  [\"284\u00d76=\", \"338\u00d76=\"],\n  [\"176\u00d73=\", \"466\u00d73=\"],\n  [\"441\u00d76=\", \"985\u00d75=\"],\n  [\"675\u00d76=\", \"975\u00d75=\"],\n  [\"346\u00d76=\", \"298\u00d73=\"],\n  [\"990\u00d78=\", \"238\u00d79=\"],\n  [\"334\u00d76=\", \"506\u00d75=\"],\n  [\"931\u00d72=\", \"701\u00d78=\"],\n  [\"772\u00d77=\", \"204\u00d79=\"],\n  [\"422\u00d78=\", \"437\u00d79=\"],\n  [\"828\u00d76=\", \"861\u00d74=\"],\n  [\"652\u00d74=\", \"985\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff:\n# - update the date line\n# - update each \"NNN\u00d7N=\" multiplication prompt in the table\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-10 Saturday\", \"2024-02-11 Sunday\"),\n    @(\"670\u00d76=\", \"392\u00d79=\"),\n    @(\"394\u00d77=\", \"174\u00d75=\"),\n    @(\"268\u00d73=\", \"424\u00d73=\"),\n    @(\"704\u00d76=\", \"802\u00d76=\"),\n    @(\"468\u00d75=\", \"514\u00d76=\"),\n    @(\"989\u00d77=\", \"904\u00d72=\"),\n    @(\"579\u00d79=\", \"518\u00d79=\"),\n    @(\"946\u00d75=\", \"139\u00d75=\"),\n    @(\"882\u00d75=\", \"431\u00d72=\"),\n    @(\"652\u00d73=\", \"299\u00d72=\"),\n    @(\"147\u00d76=\", \"333\u00d74=\"),\n    @(\"939\u00d78=\", \"422\u00d72=\"),\n    @(\"680\u00d72=\", \"115\u00d79=\"),\n    @(\"284\u00d76=\", \"338\u00d76=\"),\n    @(\"176\u00d73=\", \"466\u00d73=\"),\n    @(\"441\u00d76=\", \"985\u00d75=\"),\n    @(\"675\u00d76=\", \"975\u00d75=\"),\n    @(\"346\u00d76=\", \"298\u00d73=\"),\n    @(\"990\u00d78=\", \"238\u00d79=\"),\n    @(\"334\u00d76=\", \"506\u00d75=\"),\n    @(\"931\u00d72=\", \"701\u00d78=\"),\n    @(\"772\u00d77=\", \"204\u00d79=\"),\n    @(\"422\u00d78=\", \"437\u00d79=\"),\n    @(\"828\u00d76=\", \"861\u00d74=\"),\n    @(\"652\u00d74=\", \"985\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
